# Updated symbol list on Sat Feb  4 20:51:08 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows on Sheet1. Values are written with a leading apostrophe so
# they stay plain text (matching the source data, which stores prices and
# percentages as text) instead of being auto-coerced into numbers/percents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.24"
$ws.Range("E2").Value = "'0.22%"
$ws.Range("D3").Value = "'41.47"
$ws.Range("E3").Value = "'0.49%"
$ws.Range("D4").Value = "'5.680"
$ws.Range("E4").Value = "'-0.06%"
$ws.Range("D5").Value = "'0.08416"
$ws.Range("E5").Value = "'4.37%"
$ws.Range("D6").Value = "'8.820"
$ws.Range("E6").Value = "'0.80%"
$ws.Range("D7").Value = "'1.991"
$ws.Range("E7").Value = "'-1.38%"
$ws.Range("D8").Value = "'4.487"
$ws.Range("E8").Value = "'-0.86%"
$ws.Range("D10").Value = "'0.9255"
$ws.Range("E10").Value = "'0.05%"
$ws.Range("D11").Value = "'0.1281"
$ws.Range("E11").Value = "'1.85%"
$ws.Range("D12").Value = "'0.1964"
$ws.Range("E12").Value = "'0.94%"
$ws.Range("D13").Value = "'0.09318"
$ws.Range("E13").Value = "'0.10%"
$ws.Range("D14").Value = "'0.04001"
$ws.Range("E14").Value = "'7.99%"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("D16").Value = "'0.001305"
$ws.Range("E16").Value = "'0.39%"
$ws.Range("D17").Value = "'0.006121"
$ws.Range("E17").Value = "'-1.62%"
$ws.Range("E18").Value = "'1.77%"
$ws.Range("E19").Value = "'0.74%"
$ws.Range("D20").Value = "'8.974"
$ws.Range("E20").Value = "'8.33%"
$ws.Range("E21").Value = "'-3.85%"
$ws.Range("D22").Value = "'0.2515"
$ws.Range("E22").Value = "'-5.27%"
$ws.Range("D23").Value = "'0.04418"
$ws.Range("E23").Value = "'-0.05%"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("E24").Value = "'-1.12%"
$ws.Range("D25").Value = "'0.004364"
$ws.Range("E25").Value = "'0.36%"
$ws.Range("E26").Value = "'-4.14%"
$ws.Range("D27").Value = "'0.0004000"
$ws.Range("E27").Value = "'0.18%"
$ws.Range("D39").Value = "'0.02832"
$ws.Range("E39").Value = "'-0.49%"
$ws.Range("D40").Value = "'0.05523"
$ws.Range("E40").Value = "'1.08%"
$ws.Range("D41").Value = "'0.007916"
$ws.Range("E41").Value = "'0.23%"
$ws.Range("D42").Value = "'0.1439"
$ws.Range("E42").Value = "'1.22%"
$ws.Range("D43").Value = "'0.008979"
$ws.Range("E43").Value = "'-9.91%"
$ws.Range("D44").Value = "'0.002096"
$ws.Range("E44").Value = "'-1.83%"
$ws.Range("D45").Value = "'0.01114"
$ws.Range("E45").Value = "'-5.75%"
$ws.Range("D46").Value = "'0.00006951"
$ws.Range("E46").Value = "'2.68%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.003303"
$ws.Range("E48").Value = "'9.93%"
$ws.Range("D49").Value = "'0.002284"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.09%"
